$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "97.644.09"
$ws.Range("E2").Value = "  -1.13%  "
$ws.Range("D3").Value = "3.401.03"
$ws.Range("E3").Value = "  +2.71%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "254.27"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.31%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "652.46"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.59%  "
$ws.Range("E7").Value = "  +0.69%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.433"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +5.09%  "
$ws.Range("E9").Value = "  +6.12%  "
$ws.Range("E10").Value = "  +0.02%  "
$ws.Range("D11").Value = "3.395.85"
$ws.Range("E11").Value = "  +2.65%  "
$ws.Range("E12").Value = "  +3.68%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "41.48"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.29%  "
$ws.Range("E14").Value = "  +19.17%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000260"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.03%  "
$ws.Range("D16").Value = "97.342.70"
$ws.Range("E16").Value = "  -1.10%  "
$ws.Range("D17").Value = "4.040.97"
$ws.Range("E17").Value = "  +2.51%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "8.54"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +31.70%  "
$ws.Range("D19").Value = "3.404.82"
$ws.Range("E19").Value = "  +2.79%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.51"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +9.96%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.503"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +44.40%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.44"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.37%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.67"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +12.55%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "505.19"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.23%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000205"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.10%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.17"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +5.57%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "98.70"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +10.34%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "12.69"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.63%  "
$ws.Range("D29").Value = "3.588.25"
$ws.Range("E29").Value = "  +2.83%  "
$ws.Range("E30").Value = "  +2.54%  "
$ws.Range("E31").Value = "  +5.10%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "11.39"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.99%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.997"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.26%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.999"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.00%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.566"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +17.05%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "29.62"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.78%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.27"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +15.24%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "7.67"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +4.12%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "530.46"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +6.91%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.42"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +13.27%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.152"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.62%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "24.71"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.06%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.855"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +7.98%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.67"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -5.18%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0418"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +18.39%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.28"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.81%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "5.45"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +13.72%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.19"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +10.54%  "
$ws.Range("E50").Value = "  +11.37%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.06"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +4.31%  "
